$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-6
$ws.Cells.Item(2,2).Value = 0.13237350103577511
$ws.Cells.Item(2,3).Value = 0.75311090782945911
$ws.Cells.Item(3,2).Value = 0.1163642619345014
$ws.Cells.Item(4,2).Value = 0.091708810634297241
$ws.Cells.Item(5,2).Value = 0.070452873774489555
$ws.Cells.Item(6,2).Value = 0.066745420738513958

# Apply the same formatting as the existing data rows to the new rows first
$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A7:B11").PasteSpecial(-4122) | Out-Null

# Add new rows 7-11
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = 0.0603571035946317

$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 0.058747983545819331

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = 0.05508471830256631

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 0.052621732868468142

$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 0.048654501400396337
